$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 113, pushing existing rows 113-200 down to 114-201.
$ws.Rows.Item(113).EntireRow.Insert()

# Populate the newly inserted row 113 with the new data record.
$ws.Range("A113").Value = 4
$ws.Range("B113").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C113").Value = "Los Lagos"
$ws.Range("D113").Value = 44651
$ws.Range("E113").Value = 10
$ws.Range("F113").Value = 100112039
$ws.Range("G113").Value = "Ciboulette"
$ws.Range("H113").Value = "Sin especificar"
$ws.Range("I113").Value = "Primera"
$ws.Range("J113").Value = 80
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 3000
$ws.Range("N113").Value = "$/docena de atados"
$ws.Range("O113").Value = "Región Metropolitana"
$ws.Range("P113").Value = 1000
$ws.Range("Q113").Value = 3
$ws.Range("R113").Value = "Hortaliza"
